# Anexo3.docx: swap jornada0<->jornada4 and jornada1<->jornada3 placeholder
# indices (jornada2 is left untouched), and relocate the stray "_GoBack"
# bookmark from the end of the document to sit right after the very last
# "${jornadaN" run (the observaciones cell of the last jornada row), which
# is where it ends up once the text swap above has been applied.

$d = $word.ActiveDocument

function Replace-AllText($oldText, $newText) {
    $rng = $d.Content
    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# Drop the existing "_GoBack" bookmark first -- it will be re-created later
# at its new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# jornada0 <-> jornada4 (via a temporary placeholder so the two passes don't
# collide with one another).
Replace-AllText "jornada0" "jornadaTMPA"
Replace-AllText "jornada4" "jornada0"
Replace-AllText "jornadaTMPA" "jornada4"

# jornada1 <-> jornada3
Replace-AllText "jornada1" "jornadaTMPB"
Replace-AllText "jornada3" "jornada1"
Replace-AllText "jornadaTMPB" "jornada3"

# jornada2 is intentionally left alone.

# Locate the last "${jornada0" occurrence in the document -- that's the
# (renamed) former jornada4 "observaciones" cell, immediately followed by
# ".observaciones}".
$docEnd = $d.Content.End
$scan = $d.Range(0, $docEnd)
$lastEnd = -1
while ($true) {
    $found = $scan.Find.Execute("`${jornada0", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }
    $lastEnd = $scan.End
    $scan.Start = $scan.End
    $scan.End = $docEnd
}

if ($lastEnd -ge 0) {
    $bmRange = $d.Range($lastEnd, $lastEnd)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
